$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header label stays the same text, but sharedStrings table gains new duplicate
# entries upstream (from the new model/dataset run); re-assign to force it.
$ws.Range("A1").Value = "HK_G_acc_LG"

$values = @(
    71.739130434782609,
    71.739130434782609,
    71.739130434782609,
    72.010869565217391,
    72.010869565217391,
    72.010869565217391,
    72.010869565217391,
    72.010869565217391,
    72.010869565217391,
    72.010869565217391,
    71.739130434782609,
    72.010869565217391,
    72.010869565217391,
    72.010869565217391,
    72.010869565217391,
    72.010869565217391,
    72.010869565217391,
    72.010869565217391,
    72.010869565217391,
    72.010869565217391,
    71.739130434782609,
    72.010869565217391,
    71.739130434782609,
    71.739130434782609,
    73.369565217391312,
    72.010869565217391,
    73.369565217391312,
    71.739130434782609,
    71.467391304347828,
    71.739130434782609,
    71.739130434782609,
    72.010869565217391,
    72.010869565217391,
    72.826086956521735,
    71.467391304347828,
    71.739130434782609,
    71.467391304347828,
    73.097826086956516,
    73.369565217391312,
    72.010869565217391,
    72.010869565217391,
    71.739130434782609,
    72.010869565217391,
    72.010869565217391,
    72.010869565217391,
    72.010869565217391,
    72.010869565217391,
    72.010869565217391
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
